{"js": "// Update the cover-page copyright year from 2021 to 2022, and drop the\n// stale \"_GoBack\" bookmark left over from the previous edit session\n// (both are what the tracked diff shows changing in word/document.xml).\n\n// 1) \"\u00a9 Crown copyright 2021\" -> \"\u00a9 Crown copyright 2022\"\nconst body = context.document.body;\nconst results = body.search(\"2021\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"2022\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Remove the leftover \"_GoBack\" bookmark (start/end pair) in the\n// subtitle paragraph.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // Bookmark may already be absent; nothing further to do.\n}\n", "ps1": "# Update the cover-page copyright year from 2021 to 2022, and drop the\n# stale \"_GoBack\" bookmark left over from the previous edit session\n# (both are what the tracked diff shows changing in word/document.xml).\n\n$d = $word.ActiveDocument\n\n# 1) \"\u00a9 Crown copyright 2021\" -> \"\u00a9 Crown copyright 2022\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2021\"\n$find.MatchCase = $true\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2022\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Remove the leftover \"_GoBack\" bookmark in the subtitle paragraph.\ntry {\n  $bookmark = $d.Bookmarks.Item(\"_GoBack\")\n  $bookmark.Delete()\n} catch {\n  # Bookmark may already be absent; nothing further to do.\n}\n"}
